$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate match data (columns F:V) among rows 18, 19, 20 ---
# New row18 <- old row19, new row19 <- old row20, new row20 <- old row18
$row18 = $ws.Range("F18:V18").Value()
$row19 = $ws.Range("F19:V19").Value()
$row20 = $ws.Range("F20:V20").Value()

$ws.Range("F18:V18").Value = $row19
$ws.Range("F19:V19").Value = $row20
$ws.Range("F20:V20").Value = $row18

# --- Append new row 84 with the new match data ---
# Copy formatting from row 83 (the previous last row) for the styled columns
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)

$ws.Range("E83").Copy()
$ws.Range("E84").PasteSpecial(-4122)

$ws.Cells.Item(84, 1).Value = 83
$ws.Cells.Item(84, 2).Value = "montenegro"
$ws.Cells.Item(84, 3).Value = "prva-crnogorska-liga"
$ws.Cells.Item(84, 4).Value = "2023-2024"
$ws.Cells.Item(84, 5).Value = 45257.5625
$ws.Cells.Item(84, 6).Value = "Jezero"
$ws.Cells.Item(84, 7).Value = 4
$ws.Cells.Item(84, 8).Value = "Rudar"
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(84, 10).Value = 1.93
$ws.Cells.Item(84, 11).Value = "26/11/2023 12:44"
$ws.Cells.Item(84, 12).Value = 1.93
$ws.Cells.Item(84, 13).Value = "26/11/2023 12:44"
$ws.Cells.Item(84, 14).Value = 3.01
$ws.Cells.Item(84, 15).Value = "26/11/2023 12:44"
$ws.Cells.Item(84, 16).Value = 3.01
$ws.Cells.Item(84, 17).Value = "26/11/2023 12:44"
$ws.Cells.Item(84, 18).Value = 4.5
$ws.Cells.Item(84, 19).Value = "26/11/2023 12:44"
$ws.Cells.Item(84, 20).Value = 4.5
$ws.Cells.Item(84, 21).Value = "26/11/2023 12:44"
$ws.Cells.Item(84, 22).Value = "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/jezero-rudar/G6gwDeLA/"
